$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.892.56"
$ws.Range("E2").Value = "'  +2.66%  "
$ws.Range("D3").Value = "'1.871.07"
$ws.Range("E3").Value = "'  +0.93%  "
$ws.Range("D5").Value = "'313.75"
$ws.Range("E5").Value = "'  +1.18%  "
$ws.Range("E6").Value = "'  -0.39%  "
$ws.Range("D7").Value = "'0.4828"
$ws.Range("E7").Value = "'  +0.90%  "
$ws.Range("D8").Value = "'0.3825"
$ws.Range("E8").Value = "'  +3.56%  "
$ws.Range("D9").Value = "'0.07377"
$ws.Range("E9").Value = "'  +1.54%  "
$ws.Range("D10").Value = "'0.9393"
$ws.Range("E10").Value = "'  +0.65%  "
$ws.Range("E11").Value = "'  +5.43%  "
$ws.Range("D12").Value = "'0.07810"
$ws.Range("E12").Value = "'  +0.13%  "
$ws.Range("D13").Value = "'1.909.34"
$ws.Range("E13").Value = "'  +4.27%  "
$ws.Range("D14").Value = "'5.496"
$ws.Range("E14").Value = "'  +1.73%  "
$ws.Range("D15").Value = "'6.613"
$ws.Range("E15").Value = "'  +1.82%  "
$ws.Range("D16").Value = "'90.92"
$ws.Range("E16").Value = "'  +1.63%  "
$ws.Range("D17").Value = "'1.012"
$ws.Range("E17").Value = "'  -0.51%  "
$ws.Range("D18").Value = "'0.000008890"
$ws.Range("E18").Value = "'  +2.08%  "
$ws.Range("D20").Value = "'28.055.30"
$ws.Range("E20").Value = "'  +3.25%  "
$ws.Range("E21").Value = "'  +1.21%  "
$ws.Range("D22").Value = "'5.125"
$ws.Range("E22").Value = "'  +1.10%  "
$ws.Range("D23").Value = "'2.136.05"
$ws.Range("E23").Value = "'  +3.93%  "
$ws.Range("E24").Value = "'  +1.59%  "
$ws.Range("D25").Value = "'1.950"
$ws.Range("E25").Value = "'  +0.08%  "
$ws.Range("D26").Value = "'156.61"
$ws.Range("E26").Value = "'  +2.21%  "
$ws.Range("D27").Value = "'18.59"
$ws.Range("E27").Value = "'  +1.07%  "
$ws.Range("D28").Value = "'2.064"
$ws.Range("E28").Value = "'  +3.62%  "
$ws.Range("D29").Value = "'115.96"
$ws.Range("E29").Value = "'  +0.90%  "
$ws.Range("D30").Value = "'4.989"
$ws.Range("E30").Value = "'  +0.97%  "
$ws.Range("D31").Value = "'0.08922"
$ws.Range("E31").Value = "'  +0.36%  "
$ws.Range("D32").Value = "'3.332"
$ws.Range("E32").Value = "'  +0.90%  "
$ws.Range("D33").Value = "'1.223"
$ws.Range("E33").Value = "'  +3.03%  "
$ws.Range("D34").Value = "'0.7658"
$ws.Range("E34").Value = "'  +3.45%  "
$ws.Range("D35").Value = "'4.674"
$ws.Range("E35").Value = "'  +3.19%  "
$ws.Range("D36").Value = "'2.720"
$ws.Range("E36").Value = "'  +1.12%  "
$ws.Range("E37").Value = "'  +1.34%  "
$ws.Range("D38").Value = "'0.02047"
$ws.Range("E38").Value = "'  +3.04%  "
$ws.Range("D39").Value = "'0.5654"
$ws.Range("E39").Value = "'  +6.84%  "
$ws.Range("D40").Value = "'0.05367"
$ws.Range("E40").Value = "'  +1.64%  "
$ws.Range("D41").Value = "'2.995"
$ws.Range("E41").Value = "'  +0.53%  "
$ws.Range("D42").Value = "'7.069"
$ws.Range("E42").Value = "'  +0.15%  "
$ws.Range("D43").Value = "'8.577"
$ws.Range("E43").Value = "'  +3.08%  "
$ws.Range("E44").Value = "'  +0.47%  "
$ws.Range("B45").Value = "'Decentraland"
$ws.Range("C45").Value = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.4898"
$ws.Range("E45").Value = "'  +2.98%  "
$ws.Range("B46").Value = "'EnergySwap"
$ws.Range("C46").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.69"
$ws.Range("E46").Value = "'  +0.58%  "
$ws.Range("D47").Value = "'105.37"
$ws.Range("E47").Value = "'  +3.19%  "
$ws.Range("E48").Value = "'  -0.45%  "
$ws.Range("D49").Value = "'1.674"
$ws.Range("E49").Value = "'  +3.14%  "
$ws.Range("D50").Value = "'67.76"
$ws.Range("E50").Value = "'  +2.88%  "
$ws.Range("D51").Value = "'0.06106"
$ws.Range("E51").Value = "'  +0.67%  "
